# HIKER-M Update Pairing Matrix and Attendance List [TV]
# Adds a new date column (H) to the attendance matrix, fills attendance
# marks for each team member (mostly a checkmark, with two "x" exceptions),
# and updates the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkMark = [char]0x2713

# --- New header date in H4 (copy style from the existing date header G4) ---
$ws.Range("G4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Value = 44321

# --- Attendance marks for H5:H15 ---
# Default style (checkmark) comes from G5, which already uses the "checked"
# style (s=6) shared by most of the column.
$checkStyleSource = "G5"
# The "x" style used for the two exception rows matches the existing G10
# cell's style (s=4), which is also used elsewhere in the sheet.
$xStyleSource = "G10"

$rows = 5..15
foreach ($r in $rows) {
    $target = "H$r"
    if ($r -eq 8 -or $r -eq 14) {
        $ws.Range($xStyleSource).Copy() | Out-Null
        $ws.Range($target).PasteSpecial(-4122) | Out-Null
        $ws.Range($target).Value = "x"
    } else {
        $ws.Range($checkStyleSource).Copy() | Out-Null
        $ws.Range($target).PasteSpecial(-4122) | Out-Null
        $ws.Range($target).Value = $checkMark
    }
}

# --- Update the selected cell shown in the saved view ---
$ws.Range("I8").Select() | Out-Null

Write-Host "done"
